$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050604470900816
$ws.Range("D2").Value = 1.048474413138584
$ws.Range("E2").Value = 1.057472712464939
$ws.Range("F2").Value = 1.067958870205128
$ws.Range("I2").Value = 1.041637422615276
$ws.Range("J2").Value = 1.055636745905366
$ws.Range("K2").Value = 1.051234156994763
$ws.Range("L2").Value = 1.060207587508129
$ws.Range("M2").Value = 1.070665336545642
$ws.Range("N2").Value = 1.057135871167779
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052138308300811
$ws.Range("D3").Value = 1.049606069174046
$ws.Range("E3").Value = 1.058852907697562
$ws.Range("F3").Value = 1.06945484357063
$ws.Range("I3").Value = 1.042041416855515
$ws.Range("J3").Value = 1.056817480259296
$ws.Range("K3").Value = 1.052177038821531
$ws.Range("L3").Value = 1.061400187677058
$ws.Range("M3").Value = 1.071975505140764
$ws.Range("N3").Value = 1.058318282299923
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053129150234697
$ws.Range("D4").Value = 1.050336659446395
$ws.Range("E4").Value = 1.059744684826229
$ws.Range("F4").Value = 1.070421675341599
$ws.Range("I4").Value = 1.042300451611387
$ws.Range("J4").Value = 1.057579457470416
$ws.Range("K4").Value = 1.052784900788185
$ws.Range("L4").Value = 1.062170048398261
$ws.Range("M4").Value = 1.072821598403503
$ws.Range("N4").Value = 1.059081341606083
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053545312681335
$ws.Range("D5").Value = 1.050643405715632
$ws.Range("E5").Value = 1.060119283027848
$ws.Range("F5").Value = 1.070827860162511
$ws.Range("I5").Value = 1.042408783322482
$ws.Range("J5").Value = 1.057899311260153
$ws.Range("K5").Value = 1.053039913228749
$ws.Range("L5").Value = 1.06249326555715
$ws.Range("M5").Value = 1.073176901933544
$ws.Range("N5").Value = 1.05940164962489
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053615165705184
$ws.Range("D6").Value = 1.050694886816022
$ws.Range("E6").Value = 1.060182162027686
$ws.Range("F6").Value = 1.070896044682392
$ws.Range("I6").Value = 1.04242693954429
$ws.Range("J6").Value = 1.057952988067014
$ws.Range("K6").Value = 1.053082699841102
$ws.Range("L6").Value = 1.062547509968367
$ws.Range("M6").Value = 1.073236536040137
$ws.Range("N6").Value = 1.059455402658973
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053134712530007
$ws.Range("D7").Value = 1.05034075974638
$ws.Range("E7").Value = 1.059749691415552
$ws.Range("F7").Value = 1.070427103860039
$ws.Range("I7").Value = 1.042301901367975
$ws.Range("J7").Value = 1.057583733256458
$ws.Range("K7").Value = 1.052788310364166
$ws.Range("L7").Value = 1.062174368932786
$ws.Range("M7").Value = 1.072826347527042
$ws.Range("N7").Value = 1.059085623464232
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051123185814208
$ws.Range("D8").Value = 1.048857209848753
$ws.Range("E8").Value = 1.057939429414015
$ws.Range("F8").Value = 1.068464686553081
$ws.Range("I8").Value = 1.041774447915589
$ws.Range("J8").Value = 1.056036206099956
$ws.Range("K8").Value = 1.051553276239819
$ws.Range("L8").Value = 1.06061101521409
$ws.Range("M8").Value = 1.071108464601654
$ws.Range("N8").Value = 1.057535898641664
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047565580851959
$ws.Range("D9").Value = 1.046230006391515
$ws.Range("E9").Value = 1.054739248530686
$ws.Range("F9").Value = 1.064997426338759
$ws.Range("I9").Value = 1.040826695685043
$ws.Range("J9").Value = 1.053293386456812
$ws.Range("K9").Value = 1.049359569771282
$ws.Range("L9").Value = 1.057841884269654
$ws.Range("M9").Value = 1.068068224505341
$ws.Range("N9").Value = 1.054789183880021
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045184552260619
$ws.Range("D10").Value = 1.044469469187934
$ws.Range("E10").Value = 1.052598468925204
$ws.Range("F10").Value = 1.062679259884803
$ws.Range("I10").Value = 1.040182392768559
$ws.Range("J10").Value = 1.051453771186619
$ws.Range("K10").Value = 1.047885077785503
$ws.Range("L10").Value = 1.055985792834024
$ws.Range("M10").Value = 1.066032168612211
$ws.Range("N10").Value = 1.052946956145137
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044151218971464
$ws.Range("D11").Value = 1.043704912351548
$ws.Range("E11").Value = 1.051669655545269
$ws.Range("F11").Value = 1.061673786938309
$ws.Range("I11").Value = 1.039900409912494
$ws.Range("J11").Value = 1.050654483994146
$ws.Range("K11").Value = 1.047243685578825
$ws.Range("L11").Value = 1.055179624828038
$ws.Range("M11").Value = 1.065148252026069
$ws.Range("N11").Value = 1.052146533873171
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043767032370145
$ws.Range("D12").Value = 1.043420580007409
$ws.Range("E12").Value = 1.051324367971478
$ws.Range("F12").Value = 1.061300046529486
$ws.Range("I12").Value = 1.039795215772714
$ws.Range("J12").Value = 1.050357176565045
$ws.Range("K12").Value = 1.04700499804004
$ws.Range("L12").Value = 1.054879799733294
$ws.Range("M12").Value = 1.06481957376906
$ws.Range("N12").Value = 1.051848804233419
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043849458183593
$ws.Range("D13").Value = 1.043481585836305
$ws.Range("E13").Value = 1.051398446354782
$ws.Range("F13").Value = 1.061380227159465
$ws.Range("I13").Value = 1.039817800829285
$ws.Range("J13").Value = 1.050420968992868
$ws.Range("K13").Value = 1.047056217604477
$ws.Range("L13").Value = 1.054944130479928
$ws.Range("M13").Value = 1.064890092442986
$ws.Range("N13").Value = 1.051912687253807
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044119469402884
$ws.Range("D14").Value = 1.043681416363434
$ws.Range("E14").Value = 1.051641119827061
$ws.Range("F14").Value = 1.061642898855839
$ws.Range("I14").Value = 1.039891723794402
$ws.Range("J14").Value = 1.050629917012035
$ws.Range("K14").Value = 1.047223964724425
$ws.Range("L14").Value = 1.055154848955595
$ws.Range("M14").Value = 1.065121090612274
$ws.Range("N14").Value = 1.052121932003128
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044285784135219
$ws.Range("D15").Value = 1.043804493047723
$ws.Range("E15").Value = 1.051790600959124
$ws.Range("F15").Value = 1.061804704442857
$ws.Range("I15").Value = 1.039937210032094
$ws.Range("J15").Value = 1.050758601351302
$ws.Range("K15").Value = 1.047327259989865
$ws.Range("L15").Value = 1.055284629192917
$ws.Range("M15").Value = 1.065263369284467
$ws.Range("N15").Value = 1.052250799088917
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045253081180482
$ws.Range("D16").Value = 1.04452016271786
$ws.Range("E16").Value = 1.052660071706594
$ws.Range("F16").Value = 1.062745953419749
$ws.Range("I16").Value = 1.040201043665866
$ws.Range("J16").Value = 1.05150675924917
$ws.Range("K16").Value = 1.047927582676594
$ws.Range("L16").Value = 1.056039242901818
$ws.Range("M16").Value = 1.066090782288176
$ws.Range("N16").Value = 1.053000019456815
$ws.Range("B17").Value = 1.019999999999999
$ws.Range("C17").Value = 1.045859209475059
$ws.Range("D17").Value = 1.044968481112859
$ws.Range("E17").Value = 1.053204968870323
$ws.Range("F17").Value = 1.063335914929201
$ws.Range("I17").Value = 1.040365735375651
$ws.Range("J17").Value = 1.051975324869044
$ws.Range("K17").Value = 1.048303361188105
$ws.Range("L17").Value = 1.056511926008712
$ws.Range("M17").Value = 1.066609177562055
$ws.Range("N17").Value = 1.053469250493614
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046212529526468
$ws.Range("D18").Value = 1.045229762674002
$ws.Range("E18").Value = 1.053522621031686
$ws.Range("F18").Value = 1.063679866885929
$ws.Range("I18").Value = 1.040461508459696
$ws.Range("J18").Value = 1.052248369191538
$ws.Range("K18").Value = 1.048522264395468
$ws.Range("L18").Value = 1.056787396441926
$ws.Range("M18").Value = 1.06691132826206
$ws.Range("N18").Value = 1.053742682570365
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046332964791252
$ws.Range("D19").Value = 1.045318816665472
$ws.Range("E19").Value = 1.053630902412519
$ws.Range("F19").Value = 1.063797118299411
$ws.Range("I19").Value = 1.040494115713985
$ws.Range("J19").Value = 1.052341426027922
$ws.Range("K19").Value = 1.04859685707666
$ws.Range("L19").Value = 1.056881284663775
$ws.Range("M19").Value = 1.067014316686607
$ws.Range("N19").Value = 1.053835871558131
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04579420092271
$ws.Range("D20").Value = 1.044920403074537
$ws.Range("E20").Value = 1.053146524943128
$ws.Range("F20").Value = 1.063272634526237
$ws.Range("I20").Value = 1.040348095408772
$ws.Range("J20").Value = 1.051925079375897
$ws.Range("K20").Value = 1.048263072927139
$ws.Range("L20").Value = 1.056461236218509
$ws.Range("M20").Value = 1.066553581520429
$ws.Range("N20").Value = 1.053418933646103
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044039967865729
$ws.Range("D21").Value = 1.043622580744017
$ws.Range("E21").Value = 1.051569666469882
$ws.Range("F21").Value = 1.061565555919352
$ws.Range("I21").Value = 1.039869967858077
$ws.Range("J21").Value = 1.050568398589883
$ws.Range("K21").Value = 1.047174579737292
$ws.Range("L21").Value = 1.055092808123558
$ws.Range("M21").Value = 1.065053077198885
$ws.Range("N21").Value = 1.052060326217761
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042934920423717
$ws.Range("D22").Value = 1.042804606997888
$ws.Range("E22").Value = 1.050576580321163
$ws.Range("F22").Value = 1.06049072184323
$ws.Range("I22").Value = 1.039566726725862
$ws.Range("J22").Value = 1.04971298549122
$ws.Range("K22").Value = 1.046487617877103
$ws.Range("L22").Value = 1.054230229941573
$ws.Range("M22").Value = 1.064107608333427
$ws.Range("N22").Value = 1.05120369833438
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043520928664155
$ws.Range("D23").Value = 1.043238420339804
$ws.Range("E23").Value = 1.05110319349677
$ws.Range("F23").Value = 1.06106065934321
$ws.Range("I23").Value = 1.03972773030168
$ws.Range("J23").Value = 1.050166687564191
$ws.Range("K23").Value = 1.046852036244767
$ws.Range("L23").Value = 1.054687709388704
$ws.Range("M23").Value = 1.064609015567361
$ws.Range("N23").Value = 1.051658044716335
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045823576182591
$ws.Range("D24").Value = 1.044942128135047
$ws.Range("E24").Value = 1.053172933785699
$ws.Range("F24").Value = 1.063301228719156
$ws.Range("I24").Value = 1.040356067042987
$ws.Range("J24").Value = 1.05194778396144
$ws.Range("K24").Value = 1.048281278329747
$ws.Range("L24").Value = 1.056484141487223
$ws.Range("M24").Value = 1.066578703660056
$ws.Range("N24").Value = 1.053441670474763
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048486905974857
$ws.Range("D25").Value = 1.046910775536251
$ws.Range("E25").Value = 1.055567832054273
$ws.Range("F25").Value = 1.065894936539809
$ws.Range("I25").Value = 1.041073898391081
$ws.Range("J25").Value = 1.054004393057292
$ws.Range("K25").Value = 1.04992879138354
$ws.Range("L25").Value = 1.058559505638376
$ws.Range("M25").Value = 1.06885579558631
$ws.Range("N25").Value = 1.055501200191427
